$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "POP": fix the spelling of 'Eifel' and re-sort the park rows
# alphabetically (Eifel, Hainich, Hunsrueck, Jasmund, Kellerwald,
# Saechs_Schw, Vorpomm). Rows 6 (Jasmund) / 7 (Kellerwald) / 9 (Vorpomm)
# already sit in the correct alphabetical slot, so only rows 3, 4, 5 and 8
# need their label + Value/Std. Error/p numbers updated.
# ---------------------------------------------------------------------------
$wsPop = $wb.Worksheets.Item("POP")

$wsPop.Range("A3").Value = "ParkEifel"
$wsPop.Range("B3").Value = 0.367
$wsPop.Range("C3").Value = 0.1751
$wsPop.Range("D3").Value = 0.0361

$wsPop.Range("A4").Value = "ParkHainich"
$wsPop.Range("B4").Value = -0.1614
$wsPop.Range("C4").Value = 0.3012
$wsPop.Range("D4").Value = 0.592

$wsPop.Range("A5").Value = "ParkHunsrueck"
$wsPop.Range("B5").Value = 0.3008
$wsPop.Range("C5").Value = 0.1811
$wsPop.Range("D5").Value = 0.0967

$wsPop.Range("A8").Value = "ParkSaechs_Schw"
$wsPop.Range("B8").Value = 0.14
$wsPop.Range("C8").Value = 0.1864
$wsPop.Range("D8").Value = 0.4526

# ---------------------------------------------------------------------------
# Sheet "Pesticide": same relabeling; only the rows that actually swap
# content (old ParkHainich <-> old ParkHunsrueck slot) need their numbers
# touched, the rest already carry the right numbers for their new label.
# ---------------------------------------------------------------------------
$wsPest = $wb.Worksheets.Item("Pesticide")

$wsPest.Range("A3").Value = "ParkEifel"
$wsPest.Range("B3").Value = ""
$wsPest.Range("C3").Value = 0
$wsPest.Range("D3").Value = ""

$wsPest.Range("A4").Value = "ParkHainich"
$wsPest.Range("B4").Value = 0.4866
$wsPest.Range("C4").Value = 0.3391
$wsPest.Range("D4").Value = 0.1512

$wsPest.Range("A5").Value = "ParkHunsrueck"
$wsPest.Range("A8").Value = "ParkSaechs_Schw"
